# Cold Spell slot review - apply commit "Added many more features"
$d = $word.ActiveDocument

function Set-ParagraphText($oldText, $newText) {
    # Exact match against a whole paragraph's text (length check plus
    # prefix match) so we never touch a similar-looking substring that is
    # embedded inside unrelated prose elsewhere in the document.
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if (($t.Length -eq $oldText.Length + 1) -and ($t -like ($oldText + "*"))) {
            $p.Range.Find.Execute($oldText, $true, $true, $false, $false, $false, `
                                   $true, 1, $false, $newText, 2) | Out-Null
            return $true
        }
    }
    return $false
}

# 1) Title (Heading1 at top, and the bold repeat near the end) - both
#    occurrences share identical whole-paragraph text.
Set-ParagraphText `
    ("Play Cold Spell Slot for Free " + [char]0x2013 + " Novomatic Fantasy Theme") `
    "Play Cold Spell Slot Free - Exciting Fantasy Theme with Ice-covered Reels" | Out-Null

# run it twice: Set-ParagraphText only patches the first match it finds,
# and there are two identical paragraphs with this text.
Set-ParagraphText `
    ("Play Cold Spell Slot for Free " + [char]0x2013 + " Novomatic Fantasy Theme") `
    "Play Cold Spell Slot Free - Exciting Fantasy Theme with Ice-covered Reels" | Out-Null

# 2) Meta description (italic paragraph near the end)
Set-ParagraphText `
    "Explore a medieval realm with Cold Spell, a Novomatic online slot game with stunning ice-covered reels and exciting win potential. Play for free now." `
    "Read our review of Cold Spell, an online slot game with a fantasy theme and ice-covered reels. Play for free and enjoy big win potential." | Out-Null

# 3) "What we like" bullet list - reword the wins bullet (careful: a very
#    similar lowercase phrase is embedded in unrelated prose elsewhere, so
#    this must be scoped to the exact bullet paragraph, not document-wide).
Set-ParagraphText `
    "Wins of up to 9,000x your bet line possible" `
    "High potential for big wins with up to 9,000x your bet line" | Out-Null

# 4) "What we don't like" bullet list - reword both bullets
Set-ParagraphText `
    "High-volatility game with less frequent wins" `
    "Wins can be less frequent due to high volatility" | Out-Null

Set-ParagraphText `
    "Only 10 paylines" `
    "Limited number of paylines" | Out-Null

# 5) Insert the new "Stunning graphics..." bullet right before the
#    "Exciting fantasy theme..." bullet (same ListBullet formatting).
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Exciting fantasy theme with ice-covered reels*") {
        $p.Range.InsertParagraphBefore()
        $newPara = $d.Paragraphs.Item($i)
        $newPara.Range.Text = "Stunning graphics that bring characters to life"
        break
    }
}

# 6) Remove the old "Stunning graphics from Novomatic" bullet (it used to
#    be the last item in "What we like"; the new bullet above replaces it).
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Stunning graphics from Novomatic*") {
        $p.Range.Delete()
        break
    }
}
